# Add files via upload
# - Sheet3: update the lookup-table values in column B (rows 20,23,24,26-36)
# - Sheet1: add a new "10-nov" column (CM) mirroring the CB/CC VLOOKUP result
# - Update the active selection on Sheet1 to CN6

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- Sheet3: refresh the raw lookup values (column B, rows 20-36) -----------
$ws3.Range("B20").Value = 15.485553978162034
$ws3.Range("B23").Value = 7.9736630868682123
$ws3.Range("B24").Value = 5.8281065197545683
$ws3.Range("B26").Value = 5.7870080869847165
$ws3.Range("B27").Value = 12.96491066137621
$ws3.Range("B28").Value = 5.9050805003373839
$ws3.Range("B29").Value = 4.6023423869781102
$ws3.Range("B30").Value = 1.7856035991001411
$ws3.Range("B31").Value = 9.3892403924732761
$ws3.Range("B32").Value = 2.82068434210432
$ws3.Range("B33").Value = 17.045007135976714
$ws3.Range("B34").Value = 9.1159769433714288
$ws3.Range("B35").Value = 7.4203273211361607
$ws3.Range("B36").Value = 47.175949093667064

# --- Sheet1: add the new "10-nov" column (CM), one column past "09-nov" (CL)
$ws1.Range("CM1").Value = "10-nov"
$ws1.Range("CM1").NumberFormat = $ws1.Range("CL1").NumberFormat()

for ($r = 2; $r -le 18; $r++) {
    $cb = $ws1.Cells.Item($r, 80).Value()
    $ws1.Cells.Item($r, 91).Value = $cb
    $ws1.Cells.Item($r, 91).NumberFormat = $ws1.Cells.Item($r, 90).NumberFormat()
}

# --- Restore the active selection recorded in the workbook -------------------
$ws1.Activate()
$ws1.Range("CN6").Select()
